$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.238.88'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").Value = '1.904.12'
$ws.Range("E3").Value = '  +0.60%  '

$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").Value = '306.14'
$ws.Range("E5").Value = '  -0.33%  '

$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").Value = '0.5407'
$ws.Range("E7").Value = '  +3.42%  '

$ws.Range("D8").Value = '0.3810'
$ws.Range("E8").Value = '  +1.30%  '

$ws.Range("D9").Value = '0.07286'
$ws.Range("E9").Value = '  +0.29%  '

$ws.Range("D10").Value = '22.15'
$ws.Range("E10").Value = '  +5.05%  '

$ws.Range("D11").Value = '0.9033'
$ws.Range("E11").Value = '  +0.33%  '

$ws.Range("D12").Value = '0.08181'
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").Value = '95.82'
$ws.Range("E13").Value = '  -0.34%  '

$ws.Range("D14").Value = '5.352'
$ws.Range("E14").Value = '  +1.18%  '

$ws.Range("D15").Value = '0.9991'
$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("D16").Value = '14.89'
$ws.Range("E16").Value = '  +2.11%  '

$ws.Range("D17").Value = '0.000008647'
$ws.Range("E17").Value = '  +0.85%  '

$ws.Range("D18").Value = '0.9993'
$ws.Range("E18").Value = '  -0.25%  '

$ws.Range("D19").Value = '27.267.52'
$ws.Range("E19").Value = '  +0.46%  '

$ws.Range("D20").Value = '1.182.70'
$ws.Range("E20").Value = '  -37.76%  '

$ws.Range("E21").Value = '  -0.66%  '

$ws.Range("D22").Value = '10.82'
$ws.Range("E22").Value = '  +1.20%  '

$ws.Range("D23").Value = '6.519'
$ws.Range("E23").Value = '  +1.76%  '

$ws.Range("D24").Value = '148.50'
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").Value = '2.310'
$ws.Range("E25").Value = '  +0.99%  '

$ws.Range("D26").Value = '18.38'
$ws.Range("E26").Value = '  +1.24%  '

$ws.Range("D27").Value = '1.757'
$ws.Range("E27").Value = '  +0.92%  '

$ws.Range("D28").Value = '116.96'
$ws.Range("E28").Value = '  +1.68%  '

$ws.Range("D29").Value = '4.863'
$ws.Range("E29").Value = '  +1.53%  '

$ws.Range("D30").Value = '4.676'
$ws.Range("E30").Value = '  -3.32%  '

$ws.Range("D31").Value = '0.09224'
$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("D32").Value = '0.8332'

$ws.Range("D33").Value = '0.05078'
$ws.Range("E33").Value = '  +0.87%  '

$ws.Range("E34").Value = '  +0.77%  '

$ws.Range("D35").Value = '3.010'
$ws.Range("E35").Value = '  +1.74%  '

$ws.Range("E36").Value = '  -2.92%  '

$ws.Range("D37").Value = '2.694'
$ws.Range("E37").Value = '  +3.87%  '

$ws.Range("D38").Value = '0.5971'
$ws.Range("E38").Value = '  +4.80%  '

$ws.Range("D39").Value = '0.02004'
$ws.Range("E39").Value = '  +1.06%  '

$ws.Range("D40").Value = '1.085'
$ws.Range("E40").Value = '  +0.88%  '

$ws.Range("D41").Value = '9.276'
$ws.Range("E41").Value = '  +2.92%  '

$ws.Range("D42").Value = '6.657'
$ws.Range("E42").Value = '  +1.64%  '

$ws.Range("D43").Value = '116.35'
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("E44").Value = '  +5.64%  '

$ws.Range("D45").Value = '0.1529'
$ws.Range("E45").Value = '  +0.87%  '

$ws.Range("D46").Value = '10.20'
$ws.Range("E46").Value = '  +1.71%  '

$ws.Range("D47").Value = '0.9989'
$ws.Range("E47").Value = '  -0.29%  '

$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("D49").Value = '38.26'
$ws.Range("E49").Value = '  +0.42%  '

$ws.Range("E50").Value = '  +2.83%  '

$ws.Range("E51").Value = '  +0.15%  '
